$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Data")

# Row 6 ("sales_orders_edit") — updated end user / reseller info
$ws.Range("AG6").Value = "11,70,5"
$ws.Range("AF6").Value = "EP1234"
$ws.Range("E6").Value = "RP1234"
$ws.Range("O6").Value = "200,INGRAM MICRO TEST ACCOUNT,ATTN TOD DEBIE 1610 E SAINT ANDREW PL SANTA ANA CA"
$ws.Range("P6").Value = "50067,IRFAN MEMON,Ingram Micro Test,0,1693 Alice Ct Annapolis MD 214016511 US,TEST@TEST.IM"

# Update the on-screen selection / scroll position to match the saved view
$ws.Range("G16").Select() | Out-Null
